$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("logs")

$features = "12 features: %ascii-adp, %digit-adp, digit-adp/ascii-adp, %keyword-name, %keyword-address, %keyword-phone, bfirst-character-digit, bfirst-character-ascii, blast-character-digit, blast-character-ascii, b#ascii >= 6, b#digit >= 7"
$modelType = "Neural-Network"
$model4000 = "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 4000"

$rows = @(
    @("20160427_084407", 0.897689768976898, 0.43),
    @("20160427_094037", 0.887788778877888, 0.42),
    @("20160427_103747", 0.891089108910891, 0.42),
    @("20160427_113512", 0.904290429042904, 0.46),
    @("20160427_123307", 0.897689768976898, 0.44)
)

$r = 27
foreach ($row in $rows) {
    $time = $row[0]
    $classifyAcc = $row[1]
    $segmentAcc = $row[2]

    $ws.Cells.Item($r, 1).Value = $time
    $ws.Cells.Item($r, 2).Value = $features
    $ws.Cells.Item($r, 3).Value = $features
    $ws.Cells.Item($r, 4).Value = $features
    $ws.Cells.Item($r, 5).Value = $modelType
    $ws.Cells.Item($r, 6).Value = $model4000
    $ws.Cells.Item($r, 7).Value = $modelType
    $ws.Cells.Item($r, 8).Value = $model4000
    $ws.Cells.Item($r, 9).Value = $modelType
    $ws.Cells.Item($r, 10).Value = $model4000
    $ws.Cells.Item($r, 11).Value = $classifyAcc
    $ws.Cells.Item($r, 12).Value = $segmentAcc

    $r = $r + 1
}
